# Generate Report for Handoff
#
# Replaces the handoff-run identifiers/status across the Overview, zh-cn and
# de-de sheets:
#   - the handoff markdown file 744ee916-...md  -> f60e3c32-...md
#   - status "Ready for handoff"                -> "Handoff transform failed"
#   - the per-language "Latest Handoff File" / datetime / reason entries are
#     cleared back out to their "nothing happened yet" defaults, since the
#     handoff transform failed before a target file was produced.

$wb = $excel.ActiveWorkbook

$oldFile = "744ee916-8fe6-421b-8f6f-1069caca09ba.md"
$newFile = "f60e3c32-a450-403e-a589-e5ae29063f69.md"
$oldStatus = "Ready for handoff"
$newStatus = "Handoff transform failed"

function Replace-DisplayHyperlink {
    param($ws, $cellAddress, $newDisplay)

    $addr = '$' + ($cellAddress -replace '(\d+)', '$$$1')
    # addr now looks like $A$2 for input A2

    # Find the existing hyperlink (if any) on this cell and update its
    # display text in place (keeps the same r:id / relationship and cell
    # style). Then also push the same text into the cell value, since
    # TextToDisplay alone only rewrites the <hyperlink display="..."/>
    # attribute, not the cell's stored text.
    $match = $null
    foreach ($hl in $ws.Hyperlinks) {
        if ($hl.Range.Address() -eq $addr) {
            $match = $hl
        }
    }

    if ($match -ne $null) {
        $match.TextToDisplay = $newDisplay
    }
    $ws.Range($cellAddress).Value = $newDisplay
}

function Remove-CellHyperlink {
    param($ws, $cellAddress)

    $addr = '$' + ($cellAddress -replace '(\d+)', '$$$1')

    # As above: gather the match first, then delete it -- calling
    # Range(...).Hyperlinks.Delete() directly removes every hyperlink on the
    # whole sheet in this host, not just the one on that range.
    $match = $null
    foreach ($hl in $ws.Hyperlinks) {
        if ($hl.Range.Address() -eq $addr) {
            $match = $hl
        }
    }

    if ($match -ne $null) {
        $match.Delete()
    }
}

# ---- Overview sheet ----
$wsOverview = $wb.Worksheets.Item("Overview")
Replace-DisplayHyperlink $wsOverview "A2" $newFile
$wsOverview.Range("B2").Value = $newStatus
$wsOverview.Range("C2").Value = $newStatus

# ---- zh-cn sheet ----
$wsZh = $wb.Worksheets.Item("zh-cn")
Replace-DisplayHyperlink $wsZh "A2" $newFile
$wsZh.Range("B2").Value = $newStatus
Remove-CellHyperlink $wsZh "C2"
$wsZh.Range("C2").ClearContents()
$wsZh.Range("D2").Value = "0001-01-01 00:00:00"
$wsZh.Range("H2").Value = "Ignored"

# ---- de-de sheet ----
$wsDe = $wb.Worksheets.Item("de-de")
Replace-DisplayHyperlink $wsDe "A2" $newFile
$wsDe.Range("B2").Value = $newStatus
Remove-CellHyperlink $wsDe "C2"
$wsDe.Range("C2").ClearContents()
$wsDe.Range("D2").Value = "0001-01-01 00:00:00"
$wsDe.Range("H2").Value = "Ignored"
